$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values for rows 2-5 (columns A:AH)
$newData = @{
    2 = @(45091.50694444445, 0.861, 2.491, 0.076, 1.32, 2.861, 0.002, 1.492, 0.701, 0.38, 0.367, 3.983, 1.944, 1.138, 0.723, 1.566, 0.64, 0.961, 0, 7.973, 3.143, 0.999, 0.92, 0.309, 0.1, 0.93, 0.721, 0.494, 0.8070000000000001, 0.966, 0.167, 2.71, 3.101, 0.612)
    3 = @(45091.51388888889, 7.103, 6.019, 0.201, 14.907, 13.12, 4.979, 16.122, 8.220000000000001, 3.731, 5.366, 7.185, 6.991, 2.075, 5.368, 8.221, 4.611, 0.637, 0.041, 77.968, 15.635, 5.361, 10.128, 5.224, 0.72, 8.680999999999999, 4.568, 3.816, 4.76, 6.472, 0, 14.912, 4.061, 6.171)
    4 = @(45091.52083333334, 8.156000000000001, 6.507, 0.226, 17.28, 14.788, 5.968, 23.342, 9.553000000000001, 4.326, 6.28, 7.663, 7.767, 2.196, 6.196, 9.332000000000001, 5.268, 0.469, 0.093, 89.664, 17.913, 6.037, 11.863, 6.121, 0.827, 11.93, 5.199, 4.407, 5.415, 7.414, 0, 21.755, 4.068, 7.153)
    5 = @(45091.52777777778, 4.34, 3.51, 0.1, 9.050000000000001, 7.85, 3.06, 15.13, 4.99, 2.28, 3.23, 4.16, 4.16, 1.18, 3.23, 5.05, 2.79, 0.34, 0.01, 44, 9.630000000000001, 3.24, 6.36, 3.21, 0.44, 7.44, 2.75, 2.29, 2.89, 3.91, 0, 14.34, 2.3, 3.74)
}

foreach ($r in $newData.Keys) {
    $vals = $newData[$r]
    for ($c = 1; $c -le $vals.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $vals[$c - 1]
    }
}

# Row 6 no longer exists in the target data - remove it entirely
$ws.Rows.Item(6).Delete()

# Column width adjustments (COM ColumnWidth reports ~0.8333 narrower than the
# raw OOXML "width" attribute for this runtime's default font metrics, so we
# compensate by that offset to land on the exact target widths).
$widthOffset = 0.8333333333333334
$ws.Columns.Item(9).ColumnWidth = 7 - $widthOffset
$ws.Columns.Item(20).ColumnWidth = 8 - $widthOffset
$ws.Columns.Item(21).ColumnWidth = 8 - $widthOffset
$ws.Columns.Item(26).ColumnWidth = 7 - $widthOffset
